$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMS")

# Row 6: Change in inventories
$ws.Range("B6").Value = 25000000.0
$ws.Range("C6").Value = 28000000.0
$ws.Range("D6").Value = 59000000.0
$ws.Range("E6").Value = 45000000.0
$ws.Range("F6").Value = 6000000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 5000000.0
$ws.Range("C8").Value = 54000000.0
$ws.Range("D8").Value = 42000000.0
$ws.Range("E8").Value = 2000000.0
$ws.Range("F8").Value = -34000000.0

# Row 26: Capital Stock Change - B26 was blank inline string, now a numeric value
$ws.Range("B26").Value = 161000000.0
